# This script applies a weekly refresh of the "Fruta / hortaliza" dataset
# (commit: "Fruta / hortaliza, semanal"). The producer/category metadata in
# columns A-C, E-J stays the same; the observation-level detail in columns
# D, K-T is redistributed across rows 2-27 (row 13 is untouched) to reflect
# the new weekly pull of prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 takes the D/K-T values previously held by row 8
$ws.Range("D2").Value = 44203
$ws.Range("M2").Value = 270
$ws.Range("N2").Value = 19000
$ws.Range("O2").Value = 20000
$ws.Range("P2").Value = 19500
$ws.Range("S2").Value = 1083
$ws.Range("T2").Value = 18
$ws.Range("K2").Value = 'Carson'
$ws.Range("L2").Value = 'Tercera'
$ws.Range("Q2").Value = '$/caja 18 kilos granel'
$ws.Range("R2").Value = 'Región de O''Higgins'

# Row 3 takes the D/K-T values previously held by row 9
$ws.Range("D3").Value = 44203
$ws.Range("M3").Value = 250
$ws.Range("N3").Value = 17000
$ws.Range("O3").Value = 18000
$ws.Range("P3").Value = 17500
$ws.Range("S3").Value = 972
$ws.Range("T3").Value = 18
$ws.Range("K3").Value = 'Flavor Crest'
$ws.Range("L3").Value = 'Tercera'
$ws.Range("Q3").Value = '$/caja 18 kilos granel'
$ws.Range("R3").Value = 'Región de O''Higgins'

# Row 4 takes the D/K-T values previously held by row 10
$ws.Range("D4").Value = 44217
$ws.Range("M4").Value = 250
$ws.Range("N4").Value = 19000
$ws.Range("O4").Value = 20000
$ws.Range("P4").Value = 19500
$ws.Range("S4").Value = 1083
$ws.Range("T4").Value = 18
$ws.Range("K4").Value = 'Flavor Crest'
$ws.Range("L4").Value = 'Segunda'
$ws.Range("Q4").Value = '$/caja 18 kilos empedrada'
$ws.Range("R4").Value = 'Región de O''Higgins'

# Row 5 takes the D/K-T values previously held by row 11
$ws.Range("D5").Value = 44217
$ws.Range("M5").Value = 270
$ws.Range("N5").Value = 17000
$ws.Range("O5").Value = 18000
$ws.Range("P5").Value = 17500
$ws.Range("S5").Value = 972
$ws.Range("T5").Value = 18
$ws.Range("K5").Value = 'Royal Glory'
$ws.Range("L5").Value = 'Tercera'
$ws.Range("Q5").Value = '$/bandeja 18 kilos granel'
$ws.Range("R5").Value = 'Región de O''Higgins'

# Row 6 takes the D/K-T values previously held by row 25
$ws.Range("D6").Value = 44161
$ws.Range("M6").Value = 250
$ws.Range("N6").Value = 24000
$ws.Range("O6").Value = 25000
$ws.Range("P6").Value = 24500
$ws.Range("S6").Value = 1361
$ws.Range("T6").Value = 18
$ws.Range("K6").Value = 'Florida King'
$ws.Range("L6").Value = 'Segunda'
$ws.Range("Q6").Value = '$/caja 18 kilos granel'
$ws.Range("R6").Value = 'Región de O''Higgins'

# Row 7 takes the D/K-T values previously held by row 20
$ws.Range("D7").Value = 44257
$ws.Range("M7").Value = 300
$ws.Range("N7").Value = 19000
$ws.Range("O7").Value = 20000
$ws.Range("P7").Value = 19500
$ws.Range("S7").Value = 1083
$ws.Range("T7").Value = 18
$ws.Range("K7").Value = 'September Sweet'
$ws.Range("L7").Value = 'Segunda'
$ws.Range("Q7").Value = '$/caja 18 kilos granel'
$ws.Range("R7").Value = 'Región de O''Higgins'

# Row 8 takes the D/K-T values previously held by row 16
$ws.Range("D8").Value = 44209
$ws.Range("M8").Value = 300
$ws.Range("N8").Value = 19000
$ws.Range("O8").Value = 20000
$ws.Range("P8").Value = 19500
$ws.Range("S8").Value = 1083
$ws.Range("T8").Value = 18
$ws.Range("K8").Value = 'Carson'
$ws.Range("L8").Value = 'Segunda'
$ws.Range("Q8").Value = '$/bandeja 18 kilos granel'
$ws.Range("R8").Value = 'Región de O''Higgins'

# Row 9 takes the D/K-T values previously held by row 17
$ws.Range("D9").Value = 44209
$ws.Range("M9").Value = 300
$ws.Range("N9").Value = 18000
$ws.Range("O9").Value = 19000
$ws.Range("P9").Value = 18500
$ws.Range("S9").Value = 1028
$ws.Range("T9").Value = 18
$ws.Range("K9").Value = 'Royal Glory'
$ws.Range("L9").Value = 'Tercera'
$ws.Range("Q9").Value = '$/caja 18 kilos granel'
$ws.Range("R9").Value = 'Región de O''Higgins'

# Row 10 takes the D/K-T values previously held by row 2
$ws.Range("D10").Value = 44223
$ws.Range("M10").Value = 300
$ws.Range("N10").Value = 20000
$ws.Range("O10").Value = 21000
$ws.Range("P10").Value = 20500
$ws.Range("S10").Value = 1139
$ws.Range("T10").Value = 18
$ws.Range("K10").Value = 'Andross'
$ws.Range("L10").Value = 'Segunda'
$ws.Range("Q10").Value = '$/caja 18 kilos granel'
$ws.Range("R10").Value = 'Provincia de San Felipe de Aconcagua'

# Row 11 takes the D/K-T values previously held by row 21
$ws.Range("D11").Value = 44222
$ws.Range("M11").Value = 300
$ws.Range("N11").Value = 20000
$ws.Range("O11").Value = 21000
$ws.Range("P11").Value = 20500
$ws.Range("S11").Value = 1139
$ws.Range("T11").Value = 18
$ws.Range("K11").Value = 'Elegant Lady'
$ws.Range("L11").Value = 'Segunda'
$ws.Range("Q11").Value = '$/bandeja 18 kilos granel'
$ws.Range("R11").Value = 'Región de O''Higgins'

# Row 12 takes the D/K-T values previously held by row 3
$ws.Range("D12").Value = 44258
$ws.Range("M12").Value = 300
$ws.Range("N12").Value = 9000
$ws.Range("O12").Value = 10000
$ws.Range("P12").Value = 9500
$ws.Range("S12").Value = 950
$ws.Range("T12").Value = 10
$ws.Range("K12").Value = 'Doctor Davis'
$ws.Range("L12").Value = 'Segunda'
$ws.Range("Q12").Value = '$/bandeja 10 kilos empedrada'
$ws.Range("R12").Value = 'Provincia de San Felipe de Aconcagua'

# Row 13: unchanged
# Row 14 takes the D/K-T values previously held by row 19
$ws.Range("D14").Value = 44216
$ws.Range("M14").Value = 270
$ws.Range("N14").Value = 17000
$ws.Range("O14").Value = 18000
$ws.Range("P14").Value = 17500
$ws.Range("S14").Value = 1094
$ws.Range("T14").Value = 16
$ws.Range("K14").Value = 'Andross'
$ws.Range("L14").Value = 'Segunda'
$ws.Range("Q14").Value = '$/caja 16 kilos empedrada'
$ws.Range("R14").Value = 'Provincia de San Felipe de Aconcagua'

# Row 15 takes the D/K-T values previously held by row 26
$ws.Range("D15").Value = 44236
$ws.Range("M15").Value = 300
$ws.Range("N15").Value = 20000
$ws.Range("O15").Value = 21000
$ws.Range("P15").Value = 20500
$ws.Range("S15").Value = 1139
$ws.Range("T15").Value = 18
$ws.Range("K15").Value = 'Doctor Davis'
$ws.Range("L15").Value = 'Segunda'
$ws.Range("Q15").Value = '$/caja 18 kilos granel'
$ws.Range("R15").Value = 'Región de O''Higgins'

# Row 16 takes the D/K-T values previously held by row 27
$ws.Range("D16").Value = 44201
$ws.Range("M16").Value = 250
$ws.Range("N16").Value = 22000
$ws.Range("O16").Value = 23000
$ws.Range("P16").Value = 22500
$ws.Range("S16").Value = 1250
$ws.Range("T16").Value = 18
$ws.Range("K16").Value = 'Flavor Crest'
$ws.Range("L16").Value = 'Segunda'
$ws.Range("Q16").Value = '$/caja 18 kilos granel'
$ws.Range("R16").Value = 'Región de O''Higgins'

# Row 17 takes the D/K-T values previously held by row 14
$ws.Range("D17").Value = 44208
$ws.Range("M17").Value = 300
$ws.Range("N17").Value = 28000
$ws.Range("O17").Value = 29000
$ws.Range("P17").Value = 28500
$ws.Range("S17").Value = 1583
$ws.Range("T17").Value = 18
$ws.Range("K17").Value = 'Rich Lady'
$ws.Range("L17").Value = 'Primera'
$ws.Range("Q17").Value = '$/bandeja 18 kilos granel'
$ws.Range("R17").Value = 'Región de O''Higgins'

# Row 18 takes the D/K-T values previously held by row 22
$ws.Range("D18").Value = 44210
$ws.Range("M18").Value = 300
$ws.Range("N18").Value = 19000
$ws.Range("O18").Value = 20000
$ws.Range("P18").Value = 19500
$ws.Range("S18").Value = 1083
$ws.Range("T18").Value = 18
$ws.Range("K18").Value = 'Carson'
$ws.Range("L18").Value = 'Segunda'
$ws.Range("Q18").Value = '$/bandeja 18 kilos granel'
$ws.Range("R18").Value = 'Región de O''Higgins'

# Row 19 takes the D/K-T values previously held by row 23
$ws.Range("D19").Value = 44210
$ws.Range("M19").Value = 270
$ws.Range("N19").Value = 19000
$ws.Range("O19").Value = 20000
$ws.Range("P19").Value = 19500
$ws.Range("S19").Value = 1083
$ws.Range("T19").Value = 18
$ws.Range("K19").Value = 'Rich Lady'
$ws.Range("L19").Value = 'Segunda'
$ws.Range("Q19").Value = '$/bandeja 18 kilos granel'
$ws.Range("R19").Value = 'Región de O''Higgins'

# Row 20 takes the D/K-T values previously held by row 24
$ws.Range("D20").Value = 44210
$ws.Range("M20").Value = 300
$ws.Range("N20").Value = 19000
$ws.Range("O20").Value = 20000
$ws.Range("P20").Value = 19500
$ws.Range("S20").Value = 1083
$ws.Range("T20").Value = 18
$ws.Range("K20").Value = 'Royal Glory'
$ws.Range("L20").Value = 'Segunda'
$ws.Range("Q20").Value = '$/bandeja 18 kilos granel'
$ws.Range("R20").Value = 'Región de O''Higgins'

# Row 21 takes the D/K-T values previously held by row 5
$ws.Range("D21").Value = 44215
$ws.Range("M21").Value = 300
$ws.Range("N21").Value = 19000
$ws.Range("O21").Value = 20000
$ws.Range("P21").Value = 19500
$ws.Range("S21").Value = 1083
$ws.Range("T21").Value = 18
$ws.Range("K21").Value = 'Andross'
$ws.Range("L21").Value = 'Segunda'
$ws.Range("Q21").Value = '$/bandeja 18 kilos granel'
$ws.Range("R21").Value = 'Región de O''Higgins'

# Row 22 takes the D/K-T values previously held by row 15
$ws.Range("D22").Value = 44278
$ws.Range("M22").Value = 250
$ws.Range("N22").Value = 23000
$ws.Range("O22").Value = 24000
$ws.Range("P22").Value = 23500
$ws.Range("S22").Value = 1306
$ws.Range("T22").Value = 18
$ws.Range("K22").Value = 'Phillips Cling'
$ws.Range("L22").Value = 'Segunda'
$ws.Range("Q22").Value = '$/bandeja 18 kilos granel'
$ws.Range("R22").Value = 'Región de O''Higgins'

# Row 23 takes the D/K-T values previously held by row 18
$ws.Range("D23").Value = 44229
$ws.Range("M23").Value = 320
$ws.Range("N23").Value = 19000
$ws.Range("O23").Value = 20000
$ws.Range("P23").Value = 19500
$ws.Range("S23").Value = 1083
$ws.Range("T23").Value = 18
$ws.Range("K23").Value = 'Doctor Davis'
$ws.Range("L23").Value = 'Primera'
$ws.Range("Q23").Value = '$/bandeja 18 kilos granel'
$ws.Range("R23").Value = 'Región de O''Higgins'

# Row 24 takes the D/K-T values previously held by row 4
$ws.Range("D24").Value = 44160
$ws.Range("M24").Value = 250
$ws.Range("N24").Value = 24000
$ws.Range("O24").Value = 25000
$ws.Range("P24").Value = 24500
$ws.Range("S24").Value = 1361
$ws.Range("T24").Value = 18
$ws.Range("K24").Value = 'Early Majestic'
$ws.Range("L24").Value = 'Segunda'
$ws.Range("Q24").Value = '$/caja 18 kilos granel'
$ws.Range("R24").Value = 'Provincia de San Felipe de Aconcagua'

# Row 25 takes the D/K-T values previously held by row 6
$ws.Range("D25").Value = 44273
$ws.Range("M25").Value = 280
$ws.Range("N25").Value = 22000
$ws.Range("O25").Value = 23000
$ws.Range("P25").Value = 22500
$ws.Range("S25").Value = 1250
$ws.Range("T25").Value = 18
$ws.Range("K25").Value = 'Doctor Davis'
$ws.Range("L25").Value = 'Segunda'
$ws.Range("Q25").Value = '$/bandeja 18 kilos granel'
$ws.Range("R25").Value = 'Provincia de San Felipe de Aconcagua'

# Row 26 takes the D/K-T values previously held by row 7
$ws.Range("D26").Value = 44273
$ws.Range("M26").Value = 250
$ws.Range("N26").Value = 22000
$ws.Range("O26").Value = 23000
$ws.Range("P26").Value = 22500
$ws.Range("S26").Value = 1250
$ws.Range("T26").Value = 18
$ws.Range("K26").Value = 'Phillips Cling'
$ws.Range("L26").Value = 'Segunda'
$ws.Range("Q26").Value = '$/bandeja 18 kilos granel'
$ws.Range("R26").Value = 'Región de O''Higgins'

# Row 27 takes the D/K-T values previously held by row 12
$ws.Range("D27").Value = 44167
$ws.Range("M27").Value = 250
$ws.Range("N27").Value = 25000
$ws.Range("O27").Value = 26000
$ws.Range("P27").Value = 25500
$ws.Range("S27").Value = 1417
$ws.Range("T27").Value = 18
$ws.Range("K27").Value = 'Florida King'
$ws.Range("L27").Value = 'Segunda'
$ws.Range("Q27").Value = '$/caja 18 kilos granel'
$ws.Range("R27").Value = 'Región de O''Higgins'

